$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.09559999999999999
$ws.Range("E2").Value = -0.21
$ws.Range("G2").Value = 0.2256253508826648
$ws.Range("H2").Value = 0.2256253508826648
$ws.Range("I2").Value = 0.01266296550433535
$ws.Range("J2").Value = 0.01192967643048736
$ws.Range("K2").Value = 9.77
$ws.Range("L2").Value = 0.06094442018588984
$ws.Range("M2").Value = 17.03
$ws.Range("N2").Value = 0.03610039428498749
$ws.Range("O2").Value = 1.743091095189355
$ws.Range("P2").Value = 10.5
$ws.Range("Q2").Value = 0.02225802348751431
$ws.Range("R2").Value = 1.074718526100307
$ws.Range("S2").Value = 6.530000000000001
$ws.Range("T2").Value = 0.3834409864944217
$ws.Range("U2").Value = 432.5
$ws.Range("V2").Value = 0.916818586509518
$ws.Range("W2").Value = -0.2870690095650752
$ws.Range("X2").Value = 0.04438265700927194
$ws.Range("Y2").Value = -0.3314516665743471
$ws.Range("Z2").Value = 0.1760959448831007
$ws.Range("AA2").Value = -0.04099677312100465
$ws.Range("AB2").Value = 0.04436862383690042
$ws.Range("AC2").Value = -0.08536539695790507
$ws.Range("AD2").Value = 0.545
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.545
$ws.Range("AG2").Value = -431.955
$ws.Range("AH2").Value = 0.001153964237695459
$ws.Range("AI2").Value = 0.0004668714047689415
$ws.Range("AJ2").Value = -10.85723262536131
$ws.Range("AK2").Value = -0.5878178391361443
$ws.Range("AL2").Value = 4.873
$ws.Range("AM2").Value = 4.873
$ws.Range("AN2").Value = 0.1548295454545455
$ws.Range("AO2").Value = 0.4165811615021547
$ws.Range("AP2").Value = -122.7144886363636
$ws.Range("AQ2").Value = 0.4165811615021547

# Row 3
$ws.Range("D3").Value = -0.0218
$ws.Range("E3").Value = -0.21
$ws.Range("G3").Value = 0.2563270603504218
$ws.Range("H3").Value = 0.2563270603504218
$ws.Range("I3").Value = 0.03523685918234912
$ws.Range("J3").Value = 0.03115585559337298
$ws.Range("K3").Value = 17.3
$ws.Range("L3").Value = 0.1122647631408177
$ws.Range("M3").Value = 17.03
$ws.Range("N3").Value = 0.0367263316799655
$ws.Range("O3").Value = 0.9843930635838151
$ws.Range("P3").Value = 10.5
$ws.Range("Q3").Value = 0.0226439508302782
$ws.Range("R3").Value = 0.6069364161849711
$ws.Range("S3").Value = 6.530000000000001
$ws.Range("T3").Value = 0.3834409864944217
$ws.Range("U3").Value = 432.5
$ws.Range("V3").Value = 0.9327151175328877
$ws.Range("W3").Value = 0.01414323086984958
$ws.Range("X3").Value = 0.04439986361897712
$ws.Range("Y3").Value = -0.03025663274912754
$ws.Range("Z3").Value = 0.1768215720022949
$ws.Range("AA3").Value = 0.005509027363096701
$ws.Range("AB3").Value = 0.04437179727423409
$ws.Range("AC3").Value = -0.03886276991113739
$ws.Range("AD3").Value = 0.545
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.545
$ws.Range("AG3").Value = -431.955
$ws.Range("AH3").Value = 0.00117394910015186
$ws.Range("AI3").Value = 0.0004668714047689415
$ws.Range("AJ3").Value = -13.60702472830367
$ws.Range("AK3").Value = -0.5878178391361443
$ws.Range("AL3").Value = 0.283
$ws.Range("AM3").Value = 0.283
$ws.Range("AN3").Value = 0.07956204379562044
$ws.Range("AO3").Value = 19.18727915194346
$ws.Range("AP3").Value = -63.05912408759124
$ws.Range("AQ3").Value = 19.18727915194346

# Row 4
$ws.Range("D4").Value = 0.213
$ws.Range("G4").Value = -0.5362318840579711
$ws.Range("H4").Value = -0.5362318840579711
$ws.Range("I4").Value = -0.5475040257648953
$ws.Range("J4").Value = -0.5475040257648953
$ws.Range("K4").Value = -7.53
$ws.Range("L4").Value = -1.21256038647343
$ws.Range("W4").Value = -0.58828125
$ws.Range("X4").Value = 0.04436545039956676
$ws.Range("Y4").Value = -0.6326467003995667
$ws.Range("Z4").Value = 0.1598208770846201
$ws.Range("AA4").Value = -0.08750257360510601
$ws.Range("AB4").Value = 0.04436545039956676
$ws.Range("AC4").Value = -0.1318680240046728
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AL4").Value = 4.59
$ws.Range("AM4").Value = 4.59
$ws.Range("AN4").Value = 0
$ws.Range("AO4").Value = -0.7407407407407407
$ws.Range("AP4").Value = 0
$ws.Range("AQ4").Value = -0.7407407407407407

# Row 4: AI4 and AK4 cells were removed entirely in the diff
$ws.Range("AI4").ClearContents()
$ws.Range("AK4").ClearContents()
